$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Assign a literal string. Values that Excel would otherwise auto-convert
    # to a number (plain decimals) are forced to stay text by temporarily
    # switching the cell to a text format, then the original (default) style
    # is restored so no stray formatting is left behind.
    if ($value -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
        $range.NumberFormat = "@"
        $range.Value = $value
        $range.Style = "Normal"
    } else {
        $range.Value = $value
    }
}

Set-TextValue $ws.Range("D2") '55.704.06'
Set-TextValue $ws.Range("E2") '  -2.00%  '
Set-TextValue $ws.Range("D3") '2.342.04'
Set-TextValue $ws.Range("E3") '  -2.21%  '
Set-TextValue $ws.Range("E4") '  +0.00%  '
Set-TextValue $ws.Range("D5") '503.01'
Set-TextValue $ws.Range("E5") '  -0.34%  '
Set-TextValue $ws.Range("D6") '128.69'
Set-TextValue $ws.Range("E6") '  -2.91%  '
Set-TextValue $ws.Range("E7") '  +0.01%  '
Set-TextValue $ws.Range("D8") '0.536'
Set-TextValue $ws.Range("E8") '  -2.76%  '
Set-TextValue $ws.Range("D9") '2.348.97'
Set-TextValue $ws.Range("E9") '  -2.21%  '
Set-TextValue $ws.Range("D10") '0.0979'
Set-TextValue $ws.Range("E10") '  +0.52%  '
Set-TextValue $ws.Range("E11") '  -0.29%  '
Set-TextValue $ws.Range("D12") '4.79'
Set-TextValue $ws.Range("E12") '  +3.31%  '
Set-TextValue $ws.Range("D13") '0.319'
Set-TextValue $ws.Range("E13") '  -1.32%  '
Set-TextValue $ws.Range("D14") '2.757.82'
Set-TextValue $ws.Range("E14") '  -1.64%  '
Set-TextValue $ws.Range("D15") '55.666.41'
Set-TextValue $ws.Range("E15") '  -1.94%  '
Set-TextValue $ws.Range("D16") '21.61'
Set-TextValue $ws.Range("E16") '  -0.24%  '
Set-TextValue $ws.Range("E17") '  -0.82%  '
Set-TextValue $ws.Range("D18") '2.340.50'
Set-TextValue $ws.Range("E18") '  -1.53%  '
Set-TextValue $ws.Range("D19") '9.91'
Set-TextValue $ws.Range("E19") '  -2.99%  '
Set-TextValue $ws.Range("D20") '311.03'
Set-TextValue $ws.Range("E20") '  +0.49%  '
Set-TextValue $ws.Range("E21") '  -1.86%  '
Set-TextValue $ws.Range("E22") '  -1.14%  '
Set-TextValue $ws.Range("D23") '0.998'
Set-TextValue $ws.Range("E23") '  -0.07%  '
Set-TextValue $ws.Range("D24") '65.20'
Set-TextValue $ws.Range("E24") '  -3.04%  '
Set-TextValue $ws.Range("D25") '0.997'
Set-TextValue $ws.Range("E25") '  -0.10%  '
Set-TextValue $ws.Range("E26") '  -1.51%  '
Set-TextValue $ws.Range("E27") '  -2.90%  '
Set-TextValue $ws.Range("D28") '7.06'
Set-TextValue $ws.Range("E28") '  -4.73%  '
Set-TextValue $ws.Range("D29") '171.62'
Set-TextValue $ws.Range("E29") '  -2.50%  '
Set-TextValue $ws.Range("E30") '  -1.02%  '
Set-TextValue $ws.Range("E31") '  -3.19%  '
Set-TextValue $ws.Range("E32") '  -0.04%  '
Set-TextValue $ws.Range("D33") '5.76'
Set-TextValue $ws.Range("E33") '  -1.65%  '
Set-TextValue $ws.Range("E34") '  +0.01%  '
Set-TextValue $ws.Range("E35") '  -5.51%  '
Set-TextValue $ws.Range("D36") '17.63'
Set-TextValue $ws.Range("E36") '  -1.47%  '
Set-TextValue $ws.Range("D37") '1.17'
Set-TextValue $ws.Range("E37") '  -2.23%  '
Set-TextValue $ws.Range("E38") '  -4.58%  '
Set-TextValue $ws.Range("D39") '0.820'
Set-TextValue $ws.Range("E39") '  -1.02%  '
Set-TextValue $ws.Range("D40") '36.07'
Set-TextValue $ws.Range("E40") '  -2.01%  '
Set-TextValue $ws.Range("D41") '1.37'
Set-TextValue $ws.Range("E41") '  -4.41%  '
Set-TextValue $ws.Range("D42") '3.33'
Set-TextValue $ws.Range("E42") '  -1.34%  '
Set-TextValue $ws.Range("D43") '125.82'
Set-TextValue $ws.Range("E43") '  -4.40%  '
Set-TextValue $ws.Range("E44") '  -3.84%  '
Set-TextValue $ws.Range("D45") '0.554'
Set-TextValue $ws.Range("E45") '  -2.49%  '
Set-TextValue $ws.Range("D46") '0.0890'
Set-TextValue $ws.Range("E46") '  -2.32%  '
Set-TextValue $ws.Range("D47") '236.91'
Set-TextValue $ws.Range("E47") '  -6.05%  '
Set-TextValue $ws.Range("D48") '0.0474'
Set-TextValue $ws.Range("E48") '  -2.47%  '
Set-TextValue $ws.Range("D49") '0.0205'
Set-TextValue $ws.Range("E49") '  -2.55%  '
Set-TextValue $ws.Range("D50") '16.76'
Set-TextValue $ws.Range("E50") '  -1.93%  '
Set-TextValue $ws.Range("D51") '0.953'
Set-TextValue $ws.Range("E51") '  -0.01%  '
